# Finalized Occlusion Culling test. Prepared Object Pooling test.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New entry for 2022-06-15 (serial 44727): occlusion culling test finished.
$ws.Range("O7").Value = 44727
$ws.Range("O7").NumberFormat = $ws.Range("A7").NumberFormat
$ws.Range("P7").Value = "13.30 - 15.00"
$ws.Range("R7").Value = 1.5
$ws.Range("S7").Value = "Working on occlusion culling test."

# New entry for 2022-06-16 (serial 44728): object pooling test prepared.
$ws.Range("O8").Value = 44728
$ws.Range("O8").NumberFormat = $ws.Range("A7").NumberFormat
$ws.Range("P8").Value = "8.30 - 11.30"
$ws.Range("R8").Value = 2.5
$ws.Range("S8").Value = "Occlusion culling test. Processing data. Prepare object pooling test."

# Match the author's final cursor position.
$ws.Range("S9").Select()
